$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Buddy Hield", "SG,SF", "Golden State Warriors"),
    @("Jake LaRavia", "SF,PF", "Memphis Grizzlies"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Donovan Clingan", "C", "Portland Trail Blazers"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Tyler Herro", "PG,SG", "Miami Heat")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
